# pre alfa, se puede asignar prioridad a los ramos y a las secciones,
# se muestra todas las opciones y la de mayor peso, limitado a 6 ramos por semestre
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BASES DE DATOS now opens section 34 instead of 35.
$ws.Range("D27").Value = "31, 34"

# "Electivo Profesional" rows get distinguished names so each can carry its own
# priority / weight instead of sharing one generic label.
$ws.Range("C45").Value = "Electivo Profesional-1"
$ws.Range("C46").Value = "Electivo Profesional-2"
$ws.Range("C47").Value = "Electivo Profesional-3"
$ws.Range("C49").Value = "Electivo Profesional-4"
$ws.Range("C50").Value = "Electivo Profesional-5"
$ws.Range("C51").Value = "Electivo Profesional-6"
$ws.Range("C52").Value = "Electivo Profesional-7"
$ws.Range("C54").Value = "Electivo Profesional-8"

# Scroll/selection state left behind by the author's last save.
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("C54").Select()
